$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "1.5E-01\(\pm\)5E-05"
$ws.Range("B5").Value = "1.9E-05\(\pm\)2E-08"
$ws.Range("B8").Value = "2.3E-02\(\pm\)5E-06"
$ws.Range("C8").Value = "2.4E-02\(\pm\)8E-06"

$wb.Save()
